$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection from N12 to B21
$ws.Range("B21").Select()

# Turn the I3:K4 block of ratio formulas into a single shared formula group
# (Excel does this automatically when the same formula text is entered across
# a multi-cell range at once).
$ws.Range("I3:K4").Formula = "=F3/SUM(`$F3:`$H3)"

# Clear the stray prior_sanity value that used to live in B14
$ws.Range("B14").ClearContents()

# Update the example trial mix percentages (and let E19's entropy formula
# recalculate automatically)
$ws.Range("B19").Value = 0.4
$ws.Range("C19").Value = 0.4
$ws.Range("D19").Value = 0.2
